$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 3-5: both F and G columns move from "Not implemented" (Bad) to "Implemented" (Good)
foreach ($r in 3..5) {
    $ws.Cells.Item($r, 6).Value = "Implemented"
    $ws.Cells.Item($r, 6).Style = "Good"
    $ws.Cells.Item($r, 7).Value = "Implemented"
    $ws.Cells.Item($r, 7).Style = "Good"
}

# Rows 6-11: only the F column moves from "Not implemented" (Bad) to "Implemented" (Good)
foreach ($r in 6..11) {
    $ws.Cells.Item($r, 6).Value = "Implemented"
    $ws.Cells.Item($r, 6).Style = "Good"
}

# Update the active selection shown in the sheet view
$ws.Range("G5").Select()
